# Generate Report for handoff
# Update the "Latest Handoff Datetime" for the a995aa3c-... file on the
# zh-cn and de-de report sheets to reflect a new handoff that just occurred.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-19 07:28:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-19 07:28:56"
